# Fix Treatment from Blankinship spreadsheet
# The "R_LO" sheet (already the active sheet) had its Treatment column (D)
# mislabeled as "No treatment" for rows 17:33 (Upland stations) -- these
# should read "Upland treatment" instead. The duplicate "No treatment"
# rows further down the sheet (D50:D63, a different station set) are left
# untouched, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the treatment label for the affected range.
$ws.Range("D17:D33").Value = "Upland treatment"

# Leave the selection on the cells that were just edited, matching the
# author's on-screen state after making the change.
$ws.Range("D17:D33").Select() | Out-Null
